$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; this pushes the existing rows 36..125
# down to 37..126 (values, formats, styles all travel with them), matching
# the diff where old row36 data reappears at row37, old row37 at row38, etc.,
# and a brand-new data row lands at 36 (dimension grows to A1:R126).
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with its data.
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 45162
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112040
$ws.Range("G36").Value = "Cilantro"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 900
$ws.Range("L36").Value = 1000
$ws.Range("M36").Value = 950
$ws.Range("N36").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 475
$ws.Range("Q36").Value = 2
$ws.Range("R36").Value = "Hortaliza"

# Match the D column's date-time number format used by the rest of the
# column (style index 2 in the original file -> "YYYY-MM-DD HH:MM:SS").
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
